$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "region"
$ws.Range("B1").Value = "percent"

# Country / percent data rows (23 countries, rows 2-24)
$countries = @(
    "China", "Japan", "UK", "Germany", "Brazil", "Canada", "ANZ", "Benelux",
    "Italy", "France", "Iberia", "Mexico", "Thailand", "India", "Korea",
    "Saudi", "Gulf", "Chile", "Turkey", "Poland", "S Africa", "Indonesia", "Russia"
)
$percents = @(20, 30, 40, 99, 80, 20, 40, 50, 30, 20, 50, 60, 40, 20, 50, 10, 60, 44, 22, 66, 77, 44, 22)

for ($i = 0; $i -lt $countries.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $countries[$i]
    $ws.Cells.Item($r, 2).Value = $percents[$i]
}

# Left border down the whole A2:A24 block, with an extra top border on A2
# to close off the top of the box.
$ws.Range("A2").Borders.Item(7).LineStyle = 1
$ws.Range("A2").Borders.Item(8).LineStyle = 1
$ws.Range("A3:A24").Borders.Item(7).LineStyle = 1

$ws.Range("B25").Select()
